# Weekly update: a new week of "Betarraga" (Vega Central Mapocho de Santiago)
# price data is prepended to the existing history. Insert two fresh rows right
# above the current first data block for this series (rows 700:701) which
# pushes all the existing rows (700:732) down to (702:734), and populate the
# two new rows with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new weekly pair by inserting two blank rows; everything
# below (including styles/number formats inherited from the row being pushed
# down) shifts from 700:732 -> 702:734.
$ws.Rows("700:701").Insert()

# New row 700: "Primera" quality for the week of 2023-01-13.
$ws.Range("A700").Value = 9
$ws.Range("B700").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C700").Value = "Metropolitana"
$ws.Range("D700").Value = 44939
$ws.Range("E700").Value = 13
$ws.Range("F700").Value = 100114014
$ws.Range("G700").Value = "Betarraga"
$ws.Range("H700").Value = "Sin especificar"
$ws.Range("I700").Value = "Primera"
$ws.Range("J700").Value = 10600
$ws.Range("K700").Value = 90
$ws.Range("L700").Value = 100
$ws.Range("M700").Value = 95
$ws.Range("N700").Value = "$/unidad"
$ws.Range("O700").Value = "Región Metropolitana"
$ws.Range("P700").Value = 95
$ws.Range("Q700").Value = 1
$ws.Range("R700").Value = "Hortaliza"

# New row 701: "Segunda" quality for the same week.
$ws.Range("A701").Value = 9
$ws.Range("B701").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C701").Value = "Metropolitana"
$ws.Range("D701").Value = 44939
$ws.Range("E701").Value = 13
$ws.Range("F701").Value = 100114014
$ws.Range("G701").Value = "Betarraga"
$ws.Range("H701").Value = "Sin especificar"
$ws.Range("I701").Value = "Segunda"
$ws.Range("J701").Value = 6100
$ws.Range("K701").Value = 70
$ws.Range("L701").Value = 70
$ws.Range("M701").Value = 70
$ws.Range("N701").Value = "$/unidad"
$ws.Range("O701").Value = "Región Metropolitana"
$ws.Range("P701").Value = 70
$ws.Range("Q701").Value = 1
$ws.Range("R701").Value = "Hortaliza"
